$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally has a header row (row 1: x, y, z) followed by 20 rows
# of numeric data in A2:C21.
#
# The edit:
#   1. Inserts 6 brand-new data rows at the top of the data block (pushing
#      the existing 20 rows down from A2:C21 to A8:C27).
#   2. Appends 4 brand-new data rows after the (now shifted) data block,
#      in A28:C31.
#
# Net effect: dimension grows from A1:C21 to A1:C31.
# ---------------------------------------------------------------------------

# Step 1: shift the existing 20 rows of data (rows 2-21) down by 6 rows,
# to rows 8-27. Walk from the bottom up so we never clobber a row before
# we've read it.
for ($i = 19; $i -ge 0; $i--) {
    $srcRow = 2 + $i
    $dstRow = 8 + $i
    for ($c = 1; $c -le 3; $c++) {
        $val = $ws.Cells.Item($srcRow, $c).Value2
        $ws.Cells.Item($dstRow, $c).Value2 = $val
    }
}

# Step 2: write the 6 new rows that now occupy A2:C7.
$topData = @(
    @(-0.0410806871950626, -0.0163406450301408, -0.015118914656341),
    @(-0.022754730656743,  -0.0708603709936142, -0.0030543261673301),
    @(0.0134390350431203,  -0.0337503030896186,  0.0230601634830236),
    @(-0.0311541277915239, -0.0189368221908807,  0.0268780700862407),
    @(-0.008399397134780801, -0.0177150927484035, 0.0058032199740409),
    @(-0.0368046313524246, -0.0198531206697225,  0.0143553335219621)
)
for ($r = 0; $r -lt $topData.Length; $r++) {
    $row = $topData[$r]
    for ($c = 0; $c -lt 3; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value2 = $row[$c]
    }
}

# Step 3: write the 4 new rows appended at the bottom, A28:C31.
$bottomData = @(
    @(-0.0070249503478407, 0.098502017557621,  0.1296561509370803),
    @(-0.0024434609804302, -0.0125227374956011, -0.0360410511493682),
    @(0.016951510682702,   0.0704022198915481,  0.0198531206697225),
    @(-0.00335975876078,   -0.0123700210824608,  0.0325285755097866)
)
for ($r = 0; $r -lt $bottomData.Length; $r++) {
    $row = $bottomData[$r]
    for ($c = 0; $c -lt 3; $c++) {
        $ws.Cells.Item(28 + $r, 1 + $c).Value2 = $row[$c]
    }
}

Write-Output "UsedRange after edit: $($ws.UsedRange.Address())"
